$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.401.41'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '1.869.17'
$ws.Range("E3").Value = '  -0.46%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7034'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.20%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("E8").Value = '  -0.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3133'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.49'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.90%  '

$ws.Range("E11").Value = '  -5.02%  '

$ws.Range("D12").Value = '1.890.80'
$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '93.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.06%  '

$ws.Range("E14").Value = '  -1.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7022'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.505'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008451'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.42%  '

$ws.Range("D18").Value = '29.460.27'
$ws.Range("E18").Value = '  +0.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '252.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.65%  '

$ws.Range("D20").Value = '2.145.43'
$ws.Range("E20").Value = '  +0.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.658'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.26%  '

$ws.Range("E24").Value = '  -0.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1548'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.009'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.28%  '

$ws.Range("E27").Value = '  -0.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.80'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.60%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.506'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.35%  '

$ws.Range("E30").Value = '  -2.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.255'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.09%  '

$ws.Range("E32").Value = '  +2.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05266'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.84%  '

$ws.Range("E34").Value = '  -1.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7552'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.68%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.181'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.710'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01877'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.48%  '

$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.280.90'
$ws.Range("E39").Value = '  +0.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.774'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8970'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.94%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.020'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.54%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.41%  '

$ws.Range("E45").Value = '  -0.14%  '

$ws.Range("D46").Value = '2.042.98'
$ws.Range("E46").Value = '  +0.79%  '

$ws.Range("E47").Value = '  -4.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.804'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.619'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.5196'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4294'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.06%  '
